$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores football match odds rows. Column A holds the
# row's sequential display id and must stay untouched; columns B:AC
# hold the actual match record (id, teams, odds, results, etc.).
#
# The update re-orders several match records within the sheet (their
# B:AC payloads move to different row positions) while row A-ids stay
# sequential. Concretely:
#   - rows 101 and 102 swap their B:AC payloads
#   - rows 125 and 126 swap their B:AC payloads
#   - rows 175, 177, 178, 179 rotate their B:AC payloads:
#       175 <- 178, 178 <- 177, 177 <- 179, 179 <- 175

function Get-RowData($rowNumber) {
    return $ws.Range("B" + $rowNumber + ":AC" + $rowNumber).Value()
}

function Set-RowData($rowNumber, $data) {
    $ws.Range("B" + $rowNumber + ":AC" + $rowNumber).Value = $data
}

# --- Swap rows 101 and 102 ---
$data101 = Get-RowData 101
$data102 = Get-RowData 102
Set-RowData 101 $data102
Set-RowData 102 $data101

# --- Swap rows 125 and 126 ---
$data125 = Get-RowData 125
$data126 = Get-RowData 126
Set-RowData 125 $data126
Set-RowData 126 $data125

# --- Rotate rows 175, 177, 178, 179 ---
# Snapshot all four payloads first so none are clobbered before being read.
$data175 = Get-RowData 175
$data177 = Get-RowData 177
$data178 = Get-RowData 178
$data179 = Get-RowData 179

Set-RowData 175 $data178
Set-RowData 178 $data177
Set-RowData 177 $data179
Set-RowData 179 $data175
